$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015029178380267
$ws.Range("D2").Value = 1.020898508975491
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.013349474498112
$ws.Range("I2").Value = 1.025985017275872
$ws.Range("J2").Value = 1.020257032543634
$ws.Range("K2").Value = 1.023737433276965
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.016210892846108
$ws.Range("N2").Value = 1.010807882874443

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016261521104688
$ws.Range("D3").Value = 1.021764630230072
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.015230421777562
$ws.Range("I3").Value = 1.026172400556863
$ws.Range("J3").Value = 1.021122717544598
$ws.Range("K3").Value = 1.02440972826532
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.017893555373141
$ws.Range("N3").Value = 1.011099616869766

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017057438760935
$ws.Range("D4").Value = 1.022323710285574
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.016445631929643
$ws.Range("I4").Value = 1.026291734543892
$ws.Range("J4").Value = 1.021680951734834
$ws.Range("K4").Value = 1.024842760139441
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.018980074673626
$ws.Range("N4").Value = 1.011287593742698

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017391691191299
$ws.Range("D5").Value = 1.0225584251408
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.016956068122135
$ws.Range("I5").Value = 1.026341444912957
$ws.Range("J5").Value = 1.021915177337334
$ws.Range("K5").Value = 1.025024333619337
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.019436315000299
$ws.Range("N5").Value = 1.011366430274945

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017447793199928
$ws.Range("D6").Value = 1.022597815972305
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.017041747314966
$ws.Range("I6").Value = 1.026349764693594
$ws.Range("J6").Value = 1.021954478220663
$ws.Range("K6").Value = 1.025054792942258
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.019512888899634
$ws.Range("N6").Value = 1.011379656232332

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017061906432103
$ws.Range("D7").Value = 1.022326847820376
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.016452454107502
$ws.Range("I7").Value = 1.026292400573343
$ws.Range("J7").Value = 1.021684083255799
$ws.Range("K7").Value = 1.024845188188793
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.018986173055798
$ws.Range("N7").Value = 1.011288647901174

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.015445966297584
$ws.Range("D8").Value = 1.021191501641208
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.013985548896827
$ws.Range("I8").Value = 1.026048741526085
$ws.Range("J8").Value = 1.020549995234453
$ws.Range("K8").Value = 1.023965051824303
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.016780035285445
$ws.Range("N8").Value = 1.010906641023296

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.012586811062699
$ws.Range("D9").Value = 1.019180347249626
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.009623455046858
$ws.Range("I9").Value = 1.025604669612713
$ws.Range("J9").Value = 1.018536684697836
$ws.Range("K9").Value = 1.022398787142738
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.012874547523709
$ws.Range("N9").Value = 1.010227353585033

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.010672518092039
$ws.Range("D10").Value = 1.017832318065248
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.006704365103017
$ws.Range("I10").Value = 1.025298667504187
$ws.Range("J10").Value = 1.017184198937696
$ws.Range("K10").Value = 1.021344108701945
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.010258006146867
$ws.Range("N10").Value = 1.009770283971797

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.009841587163949
$ws.Range("D11").Value = 1.017246843468693
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.005437554726969
$ws.Range("I11").Value = 1.025163789999374
$ws.Range("J11").Value = 1.016596061128532
$ws.Range("K11").Value = 1.020884889351004
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.009121783411077
$ws.Range("N11").Value = 1.009571350417989

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.009532630648451
$ws.Range("D12").Value = 1.01702910268733
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.00496656349005
$ws.Range("I12").Value = 1.025113332141752
$ws.Range("J12").Value = 1.016377219429512
$ws.Range("K12").Value = 1.020713930098511
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.00869923724308
$ws.Range("N12").Value = 1.009497302849584

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.009598917089179
$ws.Range("D13").Value = 1.01707582107654
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.005067612960963
$ws.Range("I13").Value = 1.025124171760413
$ws.Range("J13").Value = 1.016424179047047
$ws.Range("K13").Value = 1.020750618907659
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.008789897813531
$ws.Range("N13").Value = 1.009513193335381

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.00981605512204
$ws.Range("D14").Value = 1.017228850465823
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.005398631549013
$ws.Range("I14").Value = 1.025159626454841
$ws.Range("J14").Value = 1.016577979410354
$ws.Range("K14").Value = 1.020870765675634
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.00908686595151
$ws.Range("N14").Value = 1.009565232790072

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.009949799511044
$ws.Range("D15").Value = 1.017323101108176
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.005602524016854
$ws.Range("I15").Value = 1.025181423738767
$ws.Range("J15").Value = 1.016672690231786
$ws.Range("K15").Value = 1.020944740973443
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.009269770784584
$ws.Range("N15").Value = 1.009597275455284

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.010727620898242
$ws.Range("D16").Value = 1.017871136502068
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.006788378164485
$ws.Range("I16").Value = 1.025307568698311
$ws.Range("J16").Value = 1.017223178542974
$ws.Range("K16").Value = 1.021374531824458
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.010333343884544
$ws.Range("N16").Value = 1.009783464919885

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.011214979536307
$ws.Range("D17").Value = 1.018214428424035
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.007531464807947
$ws.Range("I17").Value = 1.025386058915034
$ws.Range("J17").Value = 1.017567811695006
$ws.Range("K17").Value = 1.021643446588138
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.010999616234116
$ws.Range("N17").Value = 1.009899982581824

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.011499052241382
$ws.Range("D18").Value = 1.018414494356119
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.007964623244654
$ws.Range("I18").Value = 1.025431611714428
$ws.Range("J18").Value = 1.017768589273532
$ws.Range("K18").Value = 1.021800055521324
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.011387929751376
$ws.Range("N18").Value = 1.009967847117704

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.011595880765678
$ws.Range("D19").Value = 1.018482682869384
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.0081122735505
$ws.Range("I19").Value = 1.025447105201529
$ws.Range("J19").Value = 1.017837008497779
$ws.Range("K19").Value = 1.021853413758318
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.011520282185164
$ws.Range("N19").Value = 1.009990970555705

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.011162710823286
$ws.Range("D20").Value = 1.01817761410646
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.007451766844809
$ws.Range("I20").Value = 1.025377661376259
$ws.Range("J20").Value = 1.017530860796798
$ws.Range("K20").Value = 1.021614619895921
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.01092816388763
$ws.Range("N20").Value = 1.009887491513598

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.009752122045332
$ws.Range("D21").Value = 1.017183794592605
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.005301167149094
$ws.Range("I21").Value = 1.025149195833626
$ws.Range("J21").Value = 1.016532699599426
$ws.Range("K21").Value = 1.020835396097171
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.008999430229774
$ws.Range("N21").Value = 1.00954991275389

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.008863421627581
$ws.Range("D22").Value = 1.016557380025613
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.003946437076053
$ws.Range("I22").Value = 1.025003476540477
$ws.Range("J22").Value = 1.015902908346093
$ws.Range("K22").Value = 1.020343239352242
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.007783844139544
$ws.Range("N22").Value = 1.009336767202929

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.009334712153209
$ws.Range("D23").Value = 1.016889603421335
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.004664853754937
$ws.Range("I23").Value = 1.025080922116838
$ws.Range("J23").Value = 1.016236983603028
$ws.Range("K23").Value = 1.0206043533523
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.008428530724821
$ws.Range("N23").Value = 1.009449845210769

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.011186329407847
$ws.Range("D24").Value = 1.018194249438301
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.007487779760723
$ws.Range("I24").Value = 1.025381456570714
$ws.Range("J24").Value = 1.017547558061726
$ws.Range("K24").Value = 1.021627646192486
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.010960451089451
$ws.Range("N24").Value = 1.009893135993093

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.013327388925091
$ws.Range("D25").Value = 1.019701545049702
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.010753039481323
$ws.Range("I25").Value = 1.025721222276029
$ws.Range("J25").Value = 1.019058966257078
$ws.Range("K25").Value = 1.022805541799703
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.013886422316093
$ws.Range("N25").Value = 1.010403701911381

